$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps (coin identity + link change rows, values updated) ---
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E37").Style = "Normal"

# --- Price / Volume(1h) updates for remaining rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.153.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.382.81"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.45%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +8.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0939"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.69"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +14.26%  "
$ws.Range("E15").Value = "  +2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.744.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.377.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.145.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  +9.81%  "
$ws.Range("E20").Value = "  +2.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "279.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0933"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("E35").Value = "  +5.47%  "
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  +16.30%  "
$ws.Range("E41").Value = "  +19.80%  "
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.98"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "123.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +17.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +62.00%  "
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  +1.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.597.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.58%  "
